# Automatische test-sync: 2025-06-19 21:41:50
$wb = $excel.ActiveWorkbook

# --- 1. Append new log entry to the "Logs" sheet ---
$logs = $wb.Worksheets.Item("Logs")

$newRow = 24
$logs.Cells.Item($newRow, 1).Value = "Afmelding nieuwsbrief"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Cells.Item($newRow, 4).Value = "Afmelding / Nieuwsbrief"
$logs.Cells.Item($newRow, 6).Value = "2025-06-19 21:41:10"
$logs.Cells.Item($newRow, 7).Value = "Nee"

# --- 2. Update the category counts / order on the "Dashboard" sheet ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(5, 1).Value = "Afmelding / Nieuwsbrief"
$dash.Cells.Item(5, 2).Value = 3

$dash.Cells.Item(6, 1).Value = "IT / Technisch probleem"
$dash.Cells.Item(6, 2).Value = 2

$dash.Cells.Item(7, 1).Value = "Offerte / Prijsaanvraag"
$dash.Cells.Item(7, 2).Value = 2

# --- 3. Extend the conditional formatting ranges to include the new row ---
$catFc = $logs.Range("D2:D23").FormatConditions.Item(1)
$catFc.ModifyAppliesToRange($logs.Range("D2:D24"))

$ansFc = $logs.Range("G2:G23").FormatConditions.Item(1)
$ansFc.ModifyAppliesToRange($logs.Range("G2:G24"))
